$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-8 from 2023-09-14 (45183) to 2023-09-15 (45184)
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
